# Applies the edits described by the diff:
#  - "Ngay bao cao" (report date) table cell: 05/05/2021 -> 12/05/2021
#  - "Ngay bat dau" (start date) table cell: 28/04/2021 -> 07/05/2021
#  - "Ngay ket thuc" (end date) table cell: 05/05/2021 -> 12/05/2021
#  - "[] " keyword cell -> "[1.Function;...] " (with line breaks)
#  - "Lam bai tap khong hoc ly thuyet..." -> "Hoc ve ham va huong doi tuong trong javascript"

$d = $word.ActiveDocument

# --- Table 1: report / start / end dates -----------------------------------
$infoTable = $d.Tables.Item(1)

# Row 3: "Ngay bao cao" = 05/05/2021 -> 12/05/2021 (only the leading "05" changes)
$cellReport = $infoTable.Cell(3, 2)
$rReport = $cellReport.Range
$startReport = $rReport.Start
$d.Range($startReport, $startReport + 2).Text = "12"

# Row 4: "Ngay bat dau" = 28/04/2021 -> 07/05/2021
$cellStart = $infoTable.Cell(4, 2)
$rStart = $cellStart.Range
$startStart = $rStart.Start
# edit the month part first (back-to-front) so earlier runs aren't disturbed
$d.Range($startStart + 3, $startStart + 5).Text = "05"
$d.Range($startStart + 0, $startStart + 2).Text = "07"

# Row 5: "Ngay ket thuc" = 05/05/2021 -> 12/05/2021 (only the leading "05" changes)
$cellEnd = $infoTable.Cell(5, 2)
$rEnd = $cellEnd.Range
$startEnd = $rEnd.Start
$d.Range($startEnd, $startEnd + 2).Text = "12"

# --- Table 3: keyword summary row ------------------------------------------
$summaryTable = $d.Tables.Item(3)

# Row 2, Col 1: "[] " -> "[1.Function;...10.method] " with manual line breaks
$cellKeyword = $summaryTable.Cell(2, 1)
$rKeyword = $cellKeyword.Range
$kwStart = $rKeyword.Start
$lineBreak = [char]11
$newKeywordText = "[1.Function;2.parameter;3.argrument;4.Object;" + $lineBreak + `
    "5.prototype;6.class;7.new;8.this;9.contructor;" + $lineBreak + `
    "10.method] "
$d.Range($kwStart, $kwStart + 3).Text = $newKeywordText

# Row 2, Col 2: replace the note text
$cellNote = $summaryTable.Cell(2, 2)
$cellNote.Range.Find.Execute("Làm bài tập không học lý thuyết nên không có keyword", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Học về hàm và hướng đối tượng trong javascript", 2)
